$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate and remove the data row for "Jaren Jackson Jr." (data cleanup of player_per_game_df)
$found = $ws.Cells.Find("Jaren Jackson Jr.")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
